$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.06"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.825.99"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'313.01"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.4600"
$ws.Range("E7").Value = "  +8.10%  "
$ws.Range("D8").Value = "'0.3740"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "'0.07340"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "'0.8623"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'21.02"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.830.35"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'6.709"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'93.10"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "'5.364"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'0.000008851"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'15.03"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "27.218.63"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'5.207"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'151.85"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "'2.227"
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("D27").Value = "'18.51"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "'5.274"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("D29").Value = "'117.47"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'0.08911"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'0.7680"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").Value = "'1.198"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "'2.973"
$ws.Range("E33").Value = "  +6.08%  "
$ws.Range("D34").Value = "'4.477"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'1.106"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'0.05300"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'0.5390"
$ws.Range("E39").Value = "  +7.58%  "
$ws.Range("D40").Value = "'7.210"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").Value = "'2.884"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'0.1716"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").Value = "'0.5219"
$ws.Range("E43").Value = "  +11.39%  "
$ws.Range("D44").Value = "'8.641"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'10.69"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").Value = "'1.985"
$ws.Range("E46").Value = "  +11.12%  "
$ws.Range("D47").Value = "'106.10"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'0.06485"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").Value = "'1.684"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'0.9272"
$ws.Range("E51").Value = "  +1.87%  "
